# Update "想去人数" (interest count) figures in column F across the
# workbook's sheets, reflecting the refreshed scrape at commit 456a3b4.
#
# Sheet order (per xl/workbook.xml): 1=展览, 2=演出, 3=本地生活, 4=全部类型

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 4523
$ws.Range("F4").Value = 441
$ws.Range("F5").Value = 3628
$ws.Range("F6").Value = 1048
$ws.Range("F9").Value = 357
$ws.Range("F10").Value = 353
$ws.Range("F11").Value = 2508
$ws.Range("F12").Value = 1279
$ws.Range("F13").Value = 37
$ws.Range("F15").Value = 273
$ws.Range("F16").Value = 14
$ws.Range("F17").Value = 550
$ws.Range("F18").Value = 261
$ws.Range("F20").Value = 10357
$ws.Range("F21").Value = 6032
$ws.Range("F22").Value = 15
$ws.Range("F24").Value = 394
$ws.Range("F25").Value = 215
$ws.Range("F28").Value = 837
$ws.Range("F29").Value = 17
$ws.Range("F30").Value = 169
$ws.Range("F31").Value = 853
$ws.Range("F32").Value = 3559
$ws.Range("F35").Value = 474
$ws.Range("F36").Value = 121
$ws.Range("F37").Value = 261
$ws.Range("F38").Value = 246
$ws.Range("F39").Value = 241
$ws.Range("F40").Value = 4843
$ws.Range("F42").Value = 1126
$ws.Range("F43").Value = 164
$ws.Range("F44").Value = 169
$ws.Range("F45").Value = 94
$ws.Range("F46").Value = 484

$ws = $wb.Worksheets.Item(2)
$ws.Range("F8").Value = 28
$ws.Range("F10").Value = 23
$ws.Range("F12").Value = 133
$ws.Range("F15").Value = 3556

$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 8786
$ws.Range("F4").Value = 1625

$ws = $wb.Worksheets.Item(4)
$ws.Range("F3").Value = 1625
$ws.Range("F5").Value = 4523
$ws.Range("F7").Value = 441
$ws.Range("F8").Value = 3628
$ws.Range("F9").Value = 1048
$ws.Range("F12").Value = 353
$ws.Range("F13").Value = 2508
$ws.Range("F14").Value = 28
$ws.Range("F15").Value = 1279
$ws.Range("F17").Value = 37
$ws.Range("F18").Value = 273
$ws.Range("F19").Value = 14
$ws.Range("F20").Value = 133
$ws.Range("F21").Value = 550
$ws.Range("F22").Value = 261
$ws.Range("F24").Value = 10357
$ws.Range("F25").Value = 3556
$ws.Range("F27").Value = 15
$ws.Range("F28").Value = 394
$ws.Range("F29").Value = 215
$ws.Range("F32").Value = 837
$ws.Range("F33").Value = 17
$ws.Range("F34").Value = 169
$ws.Range("F35").Value = 853
$ws.Range("F36").Value = 3559
$ws.Range("F38").Value = 121
$ws.Range("F39").Value = 261
$ws.Range("F40").Value = 246
$ws.Range("F41").Value = 241
$ws.Range("F42").Value = 4843
$ws.Range("F44").Value = 1126
$ws.Range("F45").Value = 164
$ws.Range("F46").Value = 94
$ws.Range("F47").Value = 484
